$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) / Volume(1h) (E) columns for this run.
# Both columns hold plain display text in the source data (e.g. the
# "." thousands separator in "46.590.72" is not valid numeric syntax,
# and the % cells carry padding spaces), so force Text format before
# writing each value -- this keeps Excel from "helpfully" reparsing a
# price like "0.0790" into the number 0.079 and losing the trailing
# zero the site displayed.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.590.72"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +5.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.300.44"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.31%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.83"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.55"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +12.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.568"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.72%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.52"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +9.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.43"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.95%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.650.76"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.301.18"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.81"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.814"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "46.583.30"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +5.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.98"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +5.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0943"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.01"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.21"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.53"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +6.38%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.04%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.93"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "42.39"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +8.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.23"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.89"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +5.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.03"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.93%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +14.26%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "147.52"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.25"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +14.28%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +10.25%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.14"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +19.85%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.99"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +10.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.35"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.85%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +9.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.813.40"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.12"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +20.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.196"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +6.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "73.19"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +7.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.89"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "95.64"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.526.92"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.30%  "
